# Edits for denyhosts.docx per the target diff.
# Strategy: most hunks are "split runs with identical formatting merged
# back into a single run" (no visible text change). We reproduce that by
# running a scoped Find/Replace (old text -> same text) on the exact
# paragraph, which causes the interop layer to coalesce the run. Two
# hunks are real content edits (shape extents, a run's font, and the
# Normal style's overflowPunct flag) and are handled with direct
# property writes.

$d = $word.ActiveDocument

function Merge-ParagraphRuns {
    param(
        [int]$ParaIndex,   # 1-based Paragraphs index
        [string]$Text      # exact text to find/replace within that paragraph
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $r = $p.Range
    $r.Find.Execute($Text, $true, $false, $false, $false, $false, $true, 1, $false, $Text, 2) | Out-Null
}

# "This " + "Labtainer" + " exercise explores..." -> single run
Merge-ParagraphRuns -ParaIndex 5 -Text "This Labtainer exercise explores the use of the denyhosts utility on a SSH server to limit SSH login attempts from an IP address."

# "T" + "he lab is started..." -> single run
Merge-ParagraphRuns -ParaIndex 9 -Text "The lab is started from the labtainer working directory on your linux host, e.g., a Linux VM. From there issue the command:"

# "Key " + "file" + " #1" -> single run
Merge-ParagraphRuns -ParaIndex 22 -Text "Key file #1"

# "a" + "uth.log" -> single run
Merge-ParagraphRuns -ParaIndex 24 -Text "auth.log"

# "<smart-quote> with the password " + "hank21" -> single run
Merge-ParagraphRuns -ParaIndex 31 -Text ([char]0x201D + " with the password hank21")

# "Key " + "file " + "#2" -> single run
Merge-ParagraphRuns -ParaIndex 41 -Text "Key file #2"

# "d" + "enyhosts.conf" -> single run
Merge-ParagraphRuns -ParaIndex 43 -Text "denyhosts.conf"

# "No" + "te in particular the description and values for " -> single run
Merge-ParagraphRuns -ParaIndex 50 -Text "Note in particular the description and values for "

# "Key " + "file " + "#3" -> single run (trailing " " run is untouched)
Merge-ParagraphRuns -ParaIndex 53 -Text "Key file #3"

# ./bot.py hank run switches from "Tlwg Typo" to "Tlwg Typist" (paragraph
# mark keeps "Tlwg Typo", so trim the trailing paragraph-mark character
# out of the range before touching Font.Name).
$pBot = $d.Paragraphs.Item(68)
$rBot = $pBot.Range
$rBot.End = $rBot.End - 1
$rBot.Font.Name = "Tlwg Typist"

# Drawing / shape resize (wp:extent + a:ext both move to the new size).
$shape = $d.Shapes.Item(1)
$shape.Width = 504.9
$shape.Height = 36.9

# Normal style: overflowPunct false -> true (HangingPunctuation is the
# OM property backed by w:overflowPunct).
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $true

Write-Output "done"
